$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the b.md row ---
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-26 10:36:07"

# --- zh-cn sheet ---
# Row 2 = a.md, Row 3 = b.md
$wsZhCn.Range("C2").Value = "Ready for handoff"

$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-26 10:35:58"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e1fa9b6cdb870ce4f76ff8639e3bfd02fa518675/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/185598db536d6d345a73895ad6fb30cdfce9ce85/e2e/b.md."

$wsZhCn.Columns.Item(16).ColumnWidth = 39.15

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Ready for handoff"

$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-26 10:36:07"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e1fa9b6cdb870ce4f76ff8639e3bfd02fa518675/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/185598db536d6d345a73895ad6fb30cdfce9ce85/e2e/b.md."

$wsDeDe.Columns.Item(16).ColumnWidth = 39.15
